$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "code" (column D) and "name" (column E) values for the group columns
# (codeforiati:group-code / codeforiati:group-name) were interned into the
# shared-string table in the opposite order to before. Practically this
# means each row's column D and column E values (including the header row)
# need to swap places, and the shared string that cell B91 happened to
# reuse ("United States") moved position too, which is reflected by B91
# ending up pointing at the same slot the D/E "US" value now occupies.
#
# We reproduce the net effect by swapping D<->E for every row, and fixing
# up B91 to keep showing "United States".

$firstRow = $ws.UsedRange.Row
$lastRow = $firstRow + $ws.UsedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
}

$ws.Range("B91").Value2 = "United States"
